$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Offset" scanner/arguments column pair (N:O), mirroring the
#     existing Scanner/Arguments column pairs (B:C, D:E, F:G, H:I, J:K, L:M) ---

# Row 1: block header
$ws.Cells.Item(1, 14).Value = "Offset"

# Row 2: sub-headers
$ws.Cells.Item(2, 14).Value = "Scanner"
$ws.Cells.Item(2, 15).Value = "Arguments"

# Rows 3..58: one LinearScanner entry (with its fixed arguments string) per data row
for ($r = 3; $r -le 58; $r++) {
    $ws.Cells.Item($r, 14).Value = "LinearScanner"
    $ws.Cells.Item($r, 15).Value = "[0,0,0,0];[4095,4095,4095,4095];5"
}

# Best-effort column sizing for the two new columns (matches the other
# Scanner/Arguments-style columns being width-fitted to their content)
$ws.Columns.Item(14).AutoFit()
$ws.Columns.Item(15).AutoFit()

# Scroll back to the top and select the new N2 cell (was C59 previously)
$ws.Range("N2").Select() | Out-Null
